$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Insert 4 blank rows before row 46 (pushes existing rows 46-146 down to 50-150),
# making room for the two new products (prod9, prod10).
$ws.Range("A46:A49").EntireRow.Insert()

# Column A (translation keys) first, so new shared-string ids are allocated
# in key order: prod9Name, prod9Detail, prod10Name, prod10Detail.
$ws.Range("A46").Value = "prod9Name"
$ws.Range("A47").Value = "prod9Detail"
$ws.Range("A48").Value = "prod10Name"
$ws.Range("A49").Value = "prod10Detail"

# Row 46: prod9Name
$ws.Range("B46").Value = "Buff + Short (Men) Set"
$ws.Range("C46").Value = "เซ็ตผ้าบัฟและกางเกงวิ่ง (ชาย)"

# Row 47: prod9Detail
$ws.Range("B47").Value = "Original BUFF® Coolnet UV + Men’s 3`" Running Shorts"
$ws.Range("C47").Value = "ผ้าบัฟ รุ่น Coolnet UV + กางเกงวิ่งขาสั้นชาย 3 นิ้ว"

# Row 48: prod10Name
$ws.Range("B48").Value = "Buff + Short (Women) Set"
$ws.Range("C48").Value = "เซ็ตผ้าบัฟและกางเกงวิ่ง (หญิง)"

# Row 49: prod10Detail
$ws.Range("B49").Value = "Original BUFF® Coolnet UV + Women's 5`" Running Shorts"
$ws.Range("C49").Value = "ผ้าบัฟ รุ่น Coolnet UV + กางเกงวิ่งขาสั้นหญิง 3 นิ้ว"

# Update the view state to match the saved workbook (scrolled down a bit, new selection).
$ws.Range("C50").Select()
